$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Authorization number
$ws.Range("H6").Value = "TEST-001"

# "Otros" request-type checkbox mark -> bold red "X"
$ws.Range("H9").Value = "X"
$ws.Range("H9").Font.Bold = $true
$ws.Range("H9").Font.Color = 255

# Requester block
$ws.Range("B12").Value = "Mishell Paola Sandoval Ramirez"

$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "2025/04/03"
$ws.Cells.Item(12, 7).NumberFormat = "dd\-mmm\-yyyy"

$ws.Range("B13").Value = "TEST"

$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "2025/04/03"
$ws.Cells.Item(13, 7).NumberFormat = "dd\-mmm\-yyyy"

$ws.Range("B14").Value = "TEST"
$ws.Range("G14").Value = "TEST"

# Line items table
$ws.Range("B17").Value = "10"
$ws.Range("C17").Value = "TEST"
$ws.Range("D17").Value = "TEST"
$ws.Range("D17").WrapText = $false
$ws.Range("G17").Value = "TEST"
$ws.Range("G17").WrapText = $false

$ws.Range("B18").Value = "10"
$ws.Range("C18").Value = "TEST"
$ws.Range("D18").Value = "TEST"
$ws.Range("D18").WrapText = $false
$ws.Range("G18").Value = "TEST"
$ws.Range("G18").WrapText = $false

# Observations / provider / amount
$ws.Range("B30").Value = "Transferencia Electrónica"
$ws.Range("B31").Value = "9 - MANUEL NIETO HERNANDEZ"
$ws.Range("B32").Value = "10000"

# Requester signature block
$ws.Range("B39").Value = "Mishell Paola Sandoval Ramirez"
$ws.Range("B40").Value = "TEST"
